# Applies the 27-Jun-2024 cryptos list refresh (GitHub Actions scheduled update).
# For each changed cell: Coin/Link (B/C) are plain text, Price/Volume (D/E) are
# stored as text-formatted strings (e.g. "61.569.18", "  +1.23%  ") rather than numbers,
# exactly like the source sheet. Cells whose new value would otherwise be auto-parsed
# by Excel as a number (e.g. "0.999", "1.00", "7.80") are explicitly forced to the
# "@" (Text) number format first so the literal string is preserved on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "61.569.18"; ForceText = $false }
    @{ Cell = "E2"; Value = "  +1.23%  "; ForceText = $false }
    @{ Cell = "D3"; Value = "3.451.83"; ForceText = $false }
    @{ Cell = "E3"; Value = "  +2.19%  "; ForceText = $false }
    @{ Cell = "D4"; Value = "0.999"; ForceText = $true }
    @{ Cell = "E4"; Value = "  +0.03%  "; ForceText = $false }
    @{ Cell = "D5"; Value = "579.41"; ForceText = $true }
    @{ Cell = "E5"; Value = "  +1.33%  "; ForceText = $false }
    @{ Cell = "D6"; Value = "148.55"; ForceText = $true }
    @{ Cell = "E6"; Value = "  +9.00%  "; ForceText = $false }
    @{ Cell = "D7"; Value = "3.451.88"; ForceText = $false }
    @{ Cell = "E7"; Value = "  +2.25%  "; ForceText = $false }
    @{ Cell = "E8"; Value = "  +0.05%  "; ForceText = $false }
    @{ Cell = "D9"; Value = "0.474"; ForceText = $true }
    @{ Cell = "E9"; Value = "  +0.68%  "; ForceText = $false }
    @{ Cell = "D10"; Value = "7.80"; ForceText = $true }
    @{ Cell = "E10"; Value = "  +3.78%  "; ForceText = $false }
    @{ Cell = "E11"; Value = "  +0.78%  "; ForceText = $false }
    @{ Cell = "E12"; Value = "  +1.48%  "; ForceText = $false }
    @{ Cell = "D13"; Value = "4.045.90"; ForceText = $false }
    @{ Cell = "E13"; Value = "  +2.38%  "; ForceText = $false }
    @{ Cell = "D14"; Value = "28.08"; ForceText = $true }
    @{ Cell = "E14"; Value = "  +6.40%  "; ForceText = $false }
    @{ Cell = "E15"; Value = "  -0.68%  "; ForceText = $false }
    @{ Cell = "E16"; Value = "  +1.47%  "; ForceText = $false }
    @{ Cell = "D17"; Value = "3.447.70"; ForceText = $false }
    @{ Cell = "E17"; Value = "  +2.34%  "; ForceText = $false }
    @{ Cell = "D18"; Value = "61.676.17"; ForceText = $false }
    @{ Cell = "E18"; Value = "  +1.18%  "; ForceText = $false }
    @{ Cell = "D19"; Value = "6.34"; ForceText = $true }
    @{ Cell = "E19"; Value = "  +8.81%  "; ForceText = $false }
    @{ Cell = "D20"; Value = "14.35"; ForceText = $true }
    @{ Cell = "E20"; Value = "  +2.31%  "; ForceText = $false }
    @{ Cell = "D21"; Value = "9.44"; ForceText = $true }
    @{ Cell = "E21"; Value = "  +1.11%  "; ForceText = $false }
    @{ Cell = "D22"; Value = "386.84"; ForceText = $true }
    @{ Cell = "E22"; Value = "  +2.86%  "; ForceText = $false }
    @{ Cell = "D23"; Value = "0.571"; ForceText = $true }
    @{ Cell = "E23"; Value = "  +2.66%  "; ForceText = $false }
    @{ Cell = "D24"; Value = "3.596.47"; ForceText = $false }
    @{ Cell = "E24"; Value = "  +2.68%  "; ForceText = $false }
    @{ Cell = "D25"; Value = "72.68"; ForceText = $true }
    @{ Cell = "E25"; Value = "  +2.22%  "; ForceText = $false }
    @{ Cell = "B26"; Value = "Dai"; ForceText = $false }
    @{ Cell = "C26"; Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"; ForceText = $false }
    @{ Cell = "D26"; Value = "1.00"; ForceText = $true }
    @{ Cell = "E26"; Value = "  +0.01%  "; ForceText = $false }
    @{ Cell = "B27"; Value = "LEO"; ForceText = $false }
    @{ Cell = "C27"; Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"; ForceText = $false }
    @{ Cell = "D27"; Value = "5.77"; ForceText = $true }
    @{ Cell = "E27"; Value = "  +0.96%  "; ForceText = $false }
    @{ Cell = "E28"; Value = "  -1.60%  "; ForceText = $false }
    @{ Cell = "E29"; Value = "  +7.49%  "; ForceText = $false }
    @{ Cell = "D30"; Value = "7.84"; ForceText = $true }
    @{ Cell = "E30"; Value = "  +4.32%  "; ForceText = $false }
    @{ Cell = "D31"; Value = "1.00"; ForceText = $true }
    @{ Cell = "E31"; Value = "  +0.09%  "; ForceText = $false }
    @{ Cell = "D32"; Value = "1.51"; ForceText = $true }
    @{ Cell = "E32"; Value = "  -13.66%  "; ForceText = $false }
    @{ Cell = "E33"; Value = "  +1.26%  "; ForceText = $false }
    @{ Cell = "E34"; Value = "  +1.18%  "; ForceText = $false }
    @{ Cell = "E35"; Value = "  +0.02%  "; ForceText = $false }
    @{ Cell = "D36"; Value = "23.97"; ForceText = $true }
    @{ Cell = "E36"; Value = "  +0.99%  "; ForceText = $false }
    @{ Cell = "E37"; Value = "  +3.94%  "; ForceText = $false }
    @{ Cell = "D38"; Value = "5.23"; ForceText = $true }
    @{ Cell = "E38"; Value = "  +0.38%  "; ForceText = $false }
    @{ Cell = "D39"; Value = "1.57"; ForceText = $true }
    @{ Cell = "E39"; Value = "  +2.62%  "; ForceText = $false }
    @{ Cell = "D40"; Value = "166.60"; ForceText = $true }
    @{ Cell = "E40"; Value = "  +1.04%  "; ForceText = $false }
    @{ Cell = "E41"; Value = "  +4.87%  "; ForceText = $false }
    @{ Cell = "D42"; Value = "26.08"; ForceText = $true }
    @{ Cell = "E42"; Value = "  +9.20%  "; ForceText = $false }
    @{ Cell = "E43"; Value = "  +3.44%  "; ForceText = $false }
    @{ Cell = "B44"; Value = "Filecoin"; ForceText = $false }
    @{ Cell = "C44"; Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"; ForceText = $false }
    @{ Cell = "D44"; Value = "4.51"; ForceText = $true }
    @{ Cell = "E44"; Value = "  +2.45%  "; ForceText = $false }
    @{ Cell = "B45"; Value = "FirstDigitalUSD"; ForceText = $false }
    @{ Cell = "C45"; Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"; ForceText = $false }
    @{ Cell = "D45"; Value = "1.00"; ForceText = $true }
    @{ Cell = "E45"; Value = "  +0.06%  "; ForceText = $false }
    @{ Cell = "D46"; Value = "42.37"; ForceText = $true }
    @{ Cell = "E46"; Value = "  +2.15%  "; ForceText = $false }
    @{ Cell = "E47"; Value = "  +1.49%  "; ForceText = $false }
    @{ Cell = "D48"; Value = "2.600.51"; ForceText = $false }
    @{ Cell = "E48"; Value = "  +10.20%  "; ForceText = $false }
    @{ Cell = "E49"; Value = "  -3.63%  "; ForceText = $false }
    @{ Cell = "E50"; Value = "  +2.55%  "; ForceText = $false }
    @{ Cell = "D51"; Value = "23.31"; ForceText = $true }
    @{ Cell = "E51"; Value = "  -0.39%  "; ForceText = $false }
)

foreach ($update in $updates) {
    $range = $ws.Range($update.Cell)
    if ($update.ForceText) {
        # Prevent Excel from reinterpreting numeric-looking text as a Number,
        # which would drop formatting like trailing zeros ("1.00" -> 1).
        $range.NumberFormat = "@"
    }
    $range.Value = $update.Value
}

